$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1830
$ws.Range("F5").Value = 38
$ws.Range("F7").Value = 40
$ws.Range("F8").Value = 173
$ws.Range("F9").Value = 581
$ws.Range("F11").Value = 450
$ws.Range("F12").Value = 356
$ws.Range("F13").Value = 1387
$ws.Range("F15").Value = 1410
$ws.Range("F16").Value = 16
$ws.Range("F17").Value = 1140
$ws.Range("F18").Value = 276
$ws.Range("F19").Value = 1549
$ws.Range("F20").Value = 436
$ws.Range("F21").Value = 766
$ws.Range("F25").Value = 1205
$ws.Range("F26").Value = 310
$ws.Range("F27").Value = 33
$ws.Range("F28").Value = 785
$ws.Range("F30").Value = 1012
$ws.Range("F31").Value = 225613
$ws.Range("F32").Value = 937
$ws.Range("F36").Value = 1027
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 818
$ws.Range("F39").Value = 1555
$ws.Range("F41").Value = 23
$ws.Range("F44").Value = 762
$ws.Range("F45").Value = 105
$ws.Range("F46").Value = 28

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("E3").Value = "2024.05.04 14:20-07.28 17:30"
$ws.Range("F4").Value = 113
$ws.Range("C8").Value = "上海·THE LAST BLOSSOM~Dreamer's Band Party 3rd live BLOOM乐队毕业专场（取消）"
$ws.Range("G8").Value = "不可售"
$ws.Range("F10").Value = 144
$ws.Range("F11").Value = 1382
$ws.Range("F14").Value = 2481
$ws.Range("F15").Value = 1171
$ws.Range("F17").Value = 713
$ws.Range("F18").Value = 204
$ws.Range("F20").Value = 61
$ws.Range("F23").Value = 419
$ws.Range("F24").Value = 18
$ws.Range("F27").Value = 53950
$ws.Range("F32").Value = 230
$ws.Range("F34").Value = 49
$ws.Range("F36").Value = 9
$ws.Range("F43").Value = 30
$ws.Range("F44").Value = 30
$ws.Range("F46").Value = 111

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 859
$ws.Range("F5").Value = 238
$ws.Range("F6").Value = 2681
$ws.Range("F7").Value = 4467
$ws.Range("F10").Value = 507
$ws.Range("F11").Value = 598
$ws.Range("F12").Value = 399
$ws.Range("F13").Value = 150
$ws.Range("F14").Value = 595
$ws.Range("F15").Value = 179
$ws.Range("F16").Value = 299

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1830
$ws.Range("F3").Value = 859
$ws.Range("F4").Value = 238
$ws.Range("F6").Value = 4467
$ws.Range("F7").Value = 598
$ws.Range("F8").Value = 38
$ws.Range("F9").Value = 150
$ws.Range("F10").Value = 150
$ws.Range("F11").Value = 595
$ws.Range("F12").Value = 596
$ws.Range("F13").Value = 179
$ws.Range("F16").Value = 40
$ws.Range("F17").Value = 173
$ws.Range("F18").Value = 144
$ws.Range("F19").Value = 1382
$ws.Range("F20").Value = 581
$ws.Range("F21").Value = 450
$ws.Range("F22").Value = 356
$ws.Range("F23").Value = 2481
$ws.Range("F24").Value = 1171
$ws.Range("F25").Value = 1387
$ws.Range("F27").Value = 1410
$ws.Range("F28").Value = 1140
$ws.Range("F29").Value = 204
$ws.Range("F30").Value = 61
$ws.Range("F31").Value = 1549
$ws.Range("F32").Value = 766
$ws.Range("F34").Value = 299
$ws.Range("F35").Value = 419
$ws.Range("F36").Value = 1205
$ws.Range("F37").Value = 785
$ws.Range("F39").Value = 1012
$ws.Range("F41").Value = 937
$ws.Range("F43").Value = 818
$ws.Range("F45").Value = 1556
$ws.Range("F50").Value = 30
$ws.Range("F51").Value = 762
